$d = $word.ActiveDocument

# Paragraph 1: "Привет, как дела" -> "Привет, как" + <w:tab/> + "дела"
# (the space between "как" and "дела" becomes a tab character)
$d.Content.Find.Execute("Привет, как дела", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Привет, как^tдела", 2)

# Paragraph 3: "Супер" -> <w:tab/> + "Супер"
$d.Content.Find.Execute("Супер", $true, $false, $false, $false, $false,
                         $true, 1, $false, "^tСупер", 2)
